{"js": "const replacements = [\n  [\"2025-04-27 Sunday\", \"2025-04-28 Monday\"],\n  [\"110\u00d77=\", \"347\u00d78=\"],\n  [\"324\u00d72=\", \"294\u00d79=\"],\n  [\"187\u00d78=\", \"348\u00d79=\"],\n  [\"159\u00d77=\", \"418\u00d74=\"],\n  [\"903\u00d72=\", \"267\u00d79=\"],\n  [\"578\u00d75=\", \"601\u00d79=\"],\n  [\"828\u00d76=\", \"158\u00d79=\"],\n  [\"932\u00d79=\", \"178\u00d72=\"],\n  [\"502\u00d78=\", \"529\u00d73=\"],\n  [\"231\u00d73=\", \"987\u00d73=\"],\n  [\"199\u00d77=\", \"316\u00d74=\"],\n  [\"852\u00d74=\", \"398\u00d77=\"],\n  [\"840\u00d76=\", \"442\u00d73=\"],\n  [\"414\u00d72=\", \"879\u00d75=\"],\n  [\"625\u00d73=\", \"376\u00d74=\"],\n  [\"596\u00d73=\", \"386\u00d73=\"],\n  [\"589\u00d79=\", \"790\u00d75=\"],\n  [\"745\u00d73=\", \"625\u00d74=\"],\n  [\"715\u00d76=\", \"542\u00d76=\"],\n  [\"400\u00d72=\", \"704\u00d74=\"],\n  [\"320\u00d72=\", \"251\u00d77=\"],\n  [\"827\u00d76=\", \"562\u00d78=\"],\n  [\"303\u00d77=\", \"367\u00d78=\"],\n  [\"873\u00d76=\", \"716\u00d73=\"],\n  [\"269\u00d78=\", \"445\u00d76=\"],\n];\n\nconst body = context.document.body;\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the header date and every\n# three-digit-by-one-digit multiplication prompt in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-27 Sunday\", \"2025-04-28 Monday\"),\n    @(\"110\u00d77=\", \"347\u00d78=\"),\n    @(\"324\u00d72=\", \"294\u00d79=\"),\n    @(\"187\u00d78=\", \"348\u00d79=\"),\n    @(\"159\u00d77=\", \"418\u00d74=\"),\n    @(\"903\u00d72=\", \"267\u00d79=\"),\n    @(\"578\u00d75=\", \"601\u00d79=\"),\n    @(\"828\u00d76=\", \"158\u00d79=\"),\n    @(\"932\u00d79=\", \"178\u00d72=\"),\n    @(\"502\u00d78=\", \"529\u00d73=\"),\n    @(\"231\u00d73=\", \"987\u00d73=\"),\n    @(\"199\u00d77=\", \"316\u00d74=\"),\n    @(\"852\u00d74=\", \"398\u00d77=\"),\n    @(\"840\u00d76=\", \"442\u00d73=\"),\n    @(\"414\u00d72=\", \"879\u00d75=\"),\n    @(\"625\u00d73=\", \"376\u00d74=\"),\n    @(\"596\u00d73=\", \"386\u00d73=\"),\n    @(\"589\u00d79=\", \"790\u00d75=\"),\n    @(\"745\u00d73=\", \"625\u00d74=\"),\n    @(\"715\u00d76=\", \"542\u00d76=\"),\n    @(\"400\u00d72=\", \"704\u00d74=\"),\n    @(\"320\u00d72=\", \"251\u00d77=\"),\n    @(\"827\u00d76=\", \"562\u00d78=\"),\n    @(\"303\u00d77=\", \"367\u00d78=\"),\n    @(\"873\u00d76=\", \"716\u00d73=\"),\n    @(\"269\u00d78=\", \"445\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # MatchCase=true, MatchWholeWord=false, Wrap=FindContinue(1), Replace=ReplaceAll(2)\n    $find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
